# Apply updated crypto price/volume figures (generated from diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.929.49"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "2.300.07"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'300.00"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").Value = "'97.80"
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("E7").Value = "  +1.22%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").Value = "'35.90"
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("D11").Value = "'0.0788"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "'17.72"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "'6.78"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").Value = "2.657.99"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "2.311.18"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").Value = "'0.779"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").Value = "42.878.06"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "'12.64"
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("D20").Value = "0.0₃0909"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "'6.09"
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("D22").Value = "'68.01"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").Value = "'242.68"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").Value = "'2.14"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "'25.16"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").Value = "'166.05"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").Value = "'9.04"
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("D32").Value = "'32.82"
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'4.83"
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("D35").Value = "'5.01"
$ws.Range("E35").Value = "  -3.06%  "
$ws.Range("D36").Value = "'17.33"
$ws.Range("E36").Value = "  -3.37%  "
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").Value = "'0.0688"
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").Value = "2.008.24"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("D45").Value = "'10.19"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "'2.14"
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("D47").Value = "'17.29"
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("D49").Value = "2.524.37"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "'53.28"
$ws.Range("E50").Value = "  -3.30%  "
$ws.Range("E51").Value = "  -5.30%  "
